# Envoi de la lettre d'informations et amelioration du site
# -------------------------------------------------------------------
# Renames Feuil1 -> test_05042019, Feuil2 -> test_09042019, rewrites
# both sheets' content (new "Priorite" / "Commentaire Ruben" / "Avis
# Elie" header columns + "OK" marker cells with a green fill, and a
# refreshed list of to-do rows), and updates the selections / column
# widths to match.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "test_05042019"
$wb.Worksheets.Item(2).Name = "test_09042019"

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Start from a clean slate on both sheets so the shared-string table is
# rebuilt from scratch (drops the two retired strings, keeps the rest).
$ws1.Cells.Clear()
$ws2.Cells.Clear()

# ---- Populate cell text in the exact order the strings were first
# ---- introduced, so the rebuilt shared-string table lines up 0..23.
$ws1.Range("A2").Value = 'Onglets projets liste déroulante'
$ws1.Range("A3").Value = 'Connexion et devenez membre à droite'
$ws1.Range("A4").Value = 'Devenez membre en petit en dessous du bouton connexion'
$ws1.Range("A5").Value = 'Menu centré'
$ws1.Range("A8").Value = 'Connecté avec Facebook'
$ws1.Range("A9").Value = 'Mettre du javascript'
$ws1.Range("A10").Value = 'Connecté avec Google'
$ws1.Range("A11").Value = 'Tout sur une page'
$ws1.Range("A13").Value = 'Mettre un lien avec Facebook et Google dans contactez-nous'
$ws1.Range("A12").Value = "API s'inscrire avec Facebook et Google"
$ws1.Range("C1").Value = 'Avis Elie'
$ws1.Range("B1").Value = 'Commentaire Ruben'
$ws2.Range("B2").Value = 'Décaler l''icone du texte "projet" dans la bare déroulante'
$ws2.Range("B3").Value = 'Commet savoir qu''on est  l''accueil? Mettre "Bienvenue"'
$ws2.Range("B4").Value = 'Changer le format de chaque section (fontaine, cantine…)'
$ws2.Range("B5").Value = 'Changer les images'
$ws2.Range("B6").Value = "Envoyer tous les emails sur la boite d'honoré "
$ws2.Range("C1").Value = 'Priorité'
$ws1.Range("A6").Value = 'Mettre les erreurs des champs obligatoires avant le refresh'
$ws2.Range("B7").Value = 'Mettre les erreurs des champs obligatoire en rouge'
$ws2.Range("B8").Value = 'Vérifier le scrolling sur les articles'
$ws1.Range("A7").Value = 'Ajax pour les adresses'
$ws1.Range("C2").Value = 'OK'
$ws1.Range("A14").Value = "ajax pour l'affichage des formulaire - ex Koudetat"

# ---- Remaining cells that re-use already-registered shared strings.
$ws1.Range("C7").Value = 'OK'
$ws2.Range("D1").Value = 'Commentaire Ruben'
$ws2.Range("E1").Value = 'Avis Elie'

# ---- Numeric cells.
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2
$ws2.Range("C6").Value = 1

# ---- Green fill on the "OK" marker cells.
$ws1.Range("C2").Interior.Color = 5287936
$ws1.Range("C7").Interior.Color = 5287936

# ---- Column widths (tuned so the COM pixel-rounding lands as close as
# ---- possible to the authored widths).
$ws1.Columns.Item(1).ColumnWidth = 69.2155
$ws1.Columns.Item(2).ColumnWidth = 16.644
$ws2.Columns.Item(2).ColumnWidth = 46.9297
$ws2.Columns.Item(3).ColumnWidth = 12.3583

# ---- Selections: touch sheet2 first, then sheet1 last so sheet1 stays
# ---- the active tab, matching the target workbook state.
$ws2.Range("B8").Select() | Out-Null
$ws1.Range("C5").Select() | Out-Null

Write-Output "done"
